# [2022-07-29] Update death comparison tracker for Twitter thread (Friday update - part 2)
#
# The workbook tracks, per Dutch province, weekly death counts for a
# baseline period (cols B:M) vs. the "this year" period (cols P:AA), with
# cols AE:AP holding ROUND((thisYear-base)/base*100,2) percent-difference
# formulas. This update corrects/refines several already-entered
# "this year" (2022) observations and appends the newly available week
# (row 135, ISO week 29 of 2022) with its own data + percent-diff formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously entered 2022 weekly figures (cols P:AA) ---
# Row 109 (week 3)
$ws.Range("X109").Value = 680

# Row 111 (week 5)
$ws.Range("V111").Value = 197

# Row 112 (week 6)
$ws.Range("X112").Value = 672

# Row 113 (week 7)
$ws.Range("W113").Value = 516

# Row 114 (week 8)
$ws.Range("Z114").Value = 539

# Row 123 (week 17)
$ws.Range("Q123").Value = 147

# Row 128 (week 22)
$ws.Range("V128").Value = 195

# Row 129 (week 23)
$ws.Range("X129").Value = 591

# Row 130 (week 24)
$ws.Range("W130").Value = 447
$ws.Range("X130").Value = 607
$ws.Range("AA130").Value = 228

# Row 131 (week 25)
$ws.Range("R131").Value = 105
$ws.Range("S131").Value = 221
$ws.Range("V131").Value = 212
$ws.Range("W131").Value = 443
$ws.Range("X131").Value = 699
$ws.Range("Y131").Value = 93

# Row 132 (week 26)
$ws.Range("U132").Value = 371
$ws.Range("X132").Value = 568

# Row 133 (week 27)
$ws.Range("T133").Value = 53
$ws.Range("U133").Value = 380
$ws.Range("W133").Value = 468
$ws.Range("X133").Value = 600
$ws.Range("Z133").Value = 444
$ws.Range("AA133").Value = 207

# Row 134 (week 28) - full "this year" row revised
$ws.Range("P134").Value = 121
$ws.Range("Q134").Value = 114
$ws.Range("R134").Value = 112
$ws.Range("S134").Value = 194
$ws.Range("T134").Value = 52
$ws.Range("U134").Value = 364
$ws.Range("V134").Value = 215
$ws.Range("W134").Value = 499
$ws.Range("X134").Value = 614
$ws.Range("Y134").Value = 77
$ws.Range("Z134").Value = 439
$ws.Range("AA134").Value = 240

# --- New row 135 (week 29, 2022) ---
$ws.Range("P135").Value = 123
$ws.Range("Q135").Value = 133
$ws.Range("R135").Value = 134
$ws.Range("S135").Value = 204
$ws.Range("T135").Value = 59
$ws.Range("U135").Value = 407
$ws.Range("V135").Value = 201
$ws.Range("W135").Value = 511
$ws.Range("X135").Value = 687
$ws.Range("Y135").Value = 90
$ws.Range("Z135").Value = 453
$ws.Range("AA135").Value = 239
$ws.Range("AC135").Value = 2022
$ws.Range("AD135").Value = 29

# Percent-difference formulas for the new row, copied down from row 134's
# pattern (ROUND((thisYear-base)/base*100,2) per matched column pair).
$ws.Range("AE135").Formula = "=ROUND((P135-B135)/B135*100,2)"
$ws.Range("AF135").Formula = "=ROUND((Q135-C135)/C135*100,2)"
$ws.Range("AG135").Formula = "=ROUND((R135-D135)/D135*100,2)"
$ws.Range("AH135").Formula = "=ROUND((S135-E135)/E135*100,2)"
$ws.Range("AI135").Formula = "=ROUND((T135-F135)/F135*100,2)"
$ws.Range("AJ135").Formula = "=ROUND((U135-G135)/G135*100,2)"
$ws.Range("AK135").Formula = "=ROUND((V135-H135)/H135*100,2)"
$ws.Range("AL135").Formula = "=ROUND((W135-I135)/I135*100,2)"
$ws.Range("AM135").Formula = "=ROUND((X135-J135)/J135*100,2)"
$ws.Range("AN135").Formula = "=ROUND((Y135-K135)/K135*100,2)"
$ws.Range("AO135").Formula = "=ROUND((Z135-L135)/L135*100,2)"
$ws.Range("AP135").Formula = "=ROUND((AA135-M135)/M135*100,2)"

# --- View state: the author finished by selecting the newly-filled
# formula block on the new row before saving. ---
$ws.Range("AF135:AP135").Select()
